# Auto-update draw results: append the 2025-10-10 Pick 4 draw as a new
# row (row 24) at the bottom of the results table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 24

# Columns A, C and E hold values that Excel's automatic type detection
# would otherwise coerce to a date serial / number (e.g. "2025-10-10" or
# "251010"). Mark them as Text before writing so they land as literal
# strings (matching every other row in the table), then drop the
# explicit number format again so the cells stay styled like the rest
# of the sheet (no stray formatting left behind).
$textCols = @("A", "C", "E")
foreach ($col in $textCols) {
    $ws.Range("$col$row").NumberFormat = "@"
}

$ws.Range("A$row").Value = "2025-10-10"
$ws.Range("B$row").Value = "Pick 4"
$ws.Range("C$row").Value = "251010"
$ws.Range("D$row").Value = "6-4-8-7"
$ws.Range("E$row").Value = "2025-10-10T21:37:00.965+04:00"

foreach ($col in $textCols) {
    $ws.Range("$col$row").ClearFormats()
}
